# Update column F (dSF) values on the active worksheet per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 1
    4  = -2
    5  = -3
    6  = 3
    7  = 5
    8  = 1
    9  = 2
    10 = 5
    11 = -3
    12 = -5
    13 = 5
    14 = -1
    15 = 0
    16 = -3
    17 = -1
    18 = 0
    19 = 2
    21 = 2
    22 = -1
    23 = -2
    24 = -5
    25 = 4
    26 = -1
    27 = -2
    28 = -6
    30 = -4
    32 = -3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
